## Applies the "Add files via upload" revision to gerar_corpus_iramuteq.xlsx
##
## Summary of content changes:
##  - dic_palavras_compostas!C2:C10 (the "pc"/"PC" sigla helper column) is
##    cleared out.
##  - textos_selecionados!B1:G1 header labels are renamed from the old
##    study-field headers to the generic "Variável 1".."Variável 6".
##  - textos_selecionados row 11 (the 10th sample row) is deleted, shifting
##    the trailing instruction rows up by one (old 29/30 -> new 28/29).
##  - Active sheet / selection bookkeeping is updated: dic_palavras_compostas
##    becomes the active tab (cell D29 selected), dic_siglas keeps cell K23
##    selected, textos_selecionados keeps cell E20 selected.

$wb = $excel.ActiveWorkbook

$wsDic = $wb.Worksheets.Item("dic_palavras_compostas")
$wsSig = $wb.Worksheets.Item("dic_siglas")
$wsSel = $wb.Worksheets.Item("textos_selecionados")

# --- dic_palavras_compostas: drop the sigla helper values in column C ---
$wsDic.Range("C2:C10").ClearContents()

# --- textos_selecionados: rename the variable headers ---
$wsSel.Range("B1").Value = "Variável 1"
$wsSel.Range("C1").Value = "Variável 2"
$wsSel.Range("D1").Value = "Variável 3"
$wsSel.Range("E1").Value = "Variável 4"
$wsSel.Range("F1").Value = "Variável 5"
$wsSel.Range("G1").Value = "Variável 6"

# --- textos_selecionados: remove the extra sample row (old row 11) ---
$wsSel.Rows.Item(11).Delete()

# --- restore view/selection state: dic_palavras_compostas active, cell D29 ---
$wsSig.Activate()
$wsSig.Range("K23").Select()

$wsSel.Activate()
$wsSel.Range("E20").Select()

$wsDic.Activate()
$wsDic.Range("D29").Select()

Write-Output "edit complete"
